$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing average-based metrics for row 3 (matching pattern from row 2)
$ws.Range("Q3").Value = 0.91
$ws.Range("Z3").Value = 43208.45
$ws.Range("AG3").Value = "null"

# Update the view: scroll down one row and change the active selection
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("AG4").Select()
